# Append new log rows (216-226) to Sheet1, continuing the existing
# Date / From / To / Location / Category1 / Category2 table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A216").Value = 45498
$ws.Range("B216").Value = "00:00:00"
$ws.Range("C216").Value = "02:49:00"
$ws.Range("D216").Value = "GOWA_Riser"
$ws.Range("E216").Value = "PL"
$ws.Range("F216").Value = "Riser"

$ws.Range("A217").Value = 45498
$ws.Range("B217").Value = "02:18:00"
$ws.Range("C217").Value = "03:07:00"
$ws.Range("D217").Value = "GOWA_Riser"
$ws.Range("E217").Value = "PL"
$ws.Range("F217").Value = "Flooding"

$ws.Range("A218").Value = 45498
$ws.Range("B218").Value = "03:07:00"
$ws.Range("C218").Value = "09:20:00"
$ws.Range("D218").Value = "GOWA_Riser"
$ws.Range("E218").Value = "PL"
$ws.Range("F218").Value = "Spool Inst."

$ws.Range("A219").Value = 45498
$ws.Range("B219").Value = "09:20:00"
$ws.Range("C219").Value = "10:55:00"
$ws.Range("D219").Value = "GOWA_Riser"
$ws.Range("E219").Value = "PL"
$ws.Range("F219").Value = "Metrology"

$ws.Range("A220").Value = 45498
$ws.Range("B220").Value = "10:55:00"
$ws.Range("C220").Value = "12:19:00"
$ws.Range("D220").Value = "GOWA_Riser"
$ws.Range("E220").Value = "PL"
$ws.Range("F220").Value = "Recover P/H"

$ws.Range("A221").Value = 45498
$ws.Range("B221").Value = "12:30:00"
$ws.Range("C221").Value = "23:59:00"
$ws.Range("D221").Value = "GOWA_Riser"
$ws.Range("E221").Value = "PL"
$ws.Range("F221").Value = "Riser protector"

$ws.Range("A222").Value = 45499
$ws.Range("B222").Value = "00:00:00"
$ws.Range("C222").Value = "04:22:00"
$ws.Range("D222").Value = "GOWA_Riser"
$ws.Range("E222").Value = "PL"
$ws.Range("F222").Value = "Riser protector"

$ws.Range("A223").Value = 45499
$ws.Range("B223").Value = "04:22:00"
$ws.Range("C223").Value = "15:10:00"
$ws.Range("D223").Value = "GOWA_Riser"
$ws.Range("E223").Value = "PL"
$ws.Range("F223").Value = "Paint"

$ws.Range("A224").Value = 45499
$ws.Range("B224").Value = "15:10:00"
$ws.Range("C224").Value = "23:59:00"
$ws.Range("D224").Value = "GOWA_Riser"
$ws.Range("E224").Value = "PL"
$ws.Range("F224").Value = "Spool Inst."

$ws.Range("A225").Value = 45500
$ws.Range("B225").Value = "00:00:00"
$ws.Range("C225").Value = "00:38:00"
$ws.Range("D225").Value = "GOWA_Riser"
$ws.Range("E225").Value = "PL"
$ws.Range("F225").Value = "Spool Inst."

$ws.Range("A226").Value = 45500
$ws.Range("B226").Value = "00:38:00"
$ws.Range("C226").Value = "23:59:00"
$ws.Range("D226").Value = "GOWA_Riser"
$ws.Range("E226").Value = "PL"
$ws.Range("F226").Value = "Hydrotest"

# Keep the view pointed at the newly-added tail of the log, matching the
# workbook's selection state after the edit.
$null = $ws.Range("C230").Select()
